$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 12 (pushing the
# existing rows 12-44 down to 13-45). Fill in the new row's data.
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44487
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 300000001
$ws.Range("G12").Value = "Rabanito"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 55
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 8000
$ws.Range("N12").Value = "$/docena de paquetes"
$ws.Range("O12").Value = "Provincia de Cautín"
$ws.Range("P12").Value = 667
$ws.Range("Q12").Value = 12
$ws.Range("R12").Value = "Hortaliza"

Write-Output "row inserted"
